# ---------------------------------------------------------------------------
# Add 2022-Q4 data:
#   - Insert a new "2022-Q4" worksheet (cloned from "2022-Q3" so it keeps the
#     same layout/formatting), placed right after "总计" (i.e. before the
#     existing "2022-Q3" tab). This automatically shifts 2022-Q3 / 2022-Q2 /
#     2022-Q1 one position to the right, matching the target tab order:
#       总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1
#   - Populate the new sheet with the Q4 fund-holding table (8 funds).
#   - Update the "总计" summary sheet: shift existing Q3/Q2 rows down one
#     slot, insert the new Q4 row at the top of the data, and append the
#     Q1 row that got pushed out.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by copying "2022-Q3" (keeps formatting,
#    column widths, header styles, borders, etc. identical to the other
#    quarterly sheets) and dropping it immediately before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Helper to write one data row (index, code, name, scale, position, ratio,
# market value, rank) into the Q4 sheet at worksheet row $r.
function Set-Q4Row($sheet, $r, $idx, $code, $name, $scale, $pos, $ratio, $mval, $rank) {
    $sheet.Range("A$r").Value = $idx
    $sheet.Range("B$r").Value = "'" + $code
    $sheet.Range("C$r").Value = "'" + $name
    $sheet.Range("D$r").Value = "'" + $scale
    $sheet.Range("E$r").Value = "'" + $pos
    $sheet.Range("F$r").Value = "'" + $ratio
    $sheet.Range("G$r").Value = "'" + $mval
    $sheet.Range("H$r").Value = $rank
}

# Existing rows 2-7 get overwritten in place with the Q4 numbers.
Set-Q4Row $q4 2 0 "001210" "天弘互联网灵活配置混合A"   "7.06" "93.30" "3.40" "0.2400" 10
Set-Q4Row $q4 3 1 "009986" "天弘创新领航混合A"         "1.89" "91.84" "3.45" "0.0652" 7
Set-Q4Row $q4 4 2 "012259" "天弘鑫悦成长混合C"         "1.17" "91.71" "3.68" "0.0431" 9
Set-Q4Row $q4 5 3 "015769" "天弘低碳经济混合A"         "1.12" "86.07" "3.43" "0.0384" 9
Set-Q4Row $q4 6 4 "015770" "天弘低碳经济混合C"         "0.99" "86.07" "3.43" "0.0340" 9
Set-Q4Row $q4 7 5 "009987" "天弘创新领航混合C"         "0.47" "91.84" "3.45" "0.0162" 7

# Rows 8 and 9 are brand new -> clone formatting from row 7 first so the
# index column (s="2") / borders / bold-header styling carry over, then fill
# in the values.
$q4.Range("A7:H7").Copy() | Out-Null
$q4.Range("A8:H8").PasteSpecial(-4122) | Out-Null
$q4.Range("A9:H9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Set-Q4Row $q4 8 6 "012258" "天弘鑫悦成长混合A"         "0.25" "91.71" "3.68" "0.0092" 9
Set-Q4Row $q4 9 7 "015461" "天弘互联网灵活配置混合C"   "0.12" "93.30" "3.40" "0.0041" 10

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet.
# ---------------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Clone formatting of row 4 into the brand-new row 5 first (so A5 keeps the
# bold/bordered index style), then shift the quarter rows down:
$zj.Range("A4:D4").Copy() | Out-Null
$zj.Range("A5:D5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# row 5 <- old row 4 (2022-Q1, unchanged values, just renumbered index)
$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "'2022-Q1"
$zj.Range("C5").Value = 3
$zj.Range("D5").Value = 0.45

# row 4 <- old row 3 (2022-Q2, unchanged values)
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "'2022-Q2"
$zj.Range("C4").Value = 6
$zj.Range("D4").Value = 0.66

# row 3 <- old row 2 (2022-Q3, unchanged values)
$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "'2022-Q3"
$zj.Range("C3").Value = 6
$zj.Range("D3").Value = 0.25

# row 2 <- brand new 2022-Q4 summary row
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "'2022-Q4"
$zj.Range("C2").Value = 8
$zj.Range("D2").Value = 0.45

# Leave the original selection/active sheet pointed at 总计, like the source
# workbook.
$zj.Activate()
